# #3473 replaced two properties that had gaps
#
# This updates the "BPS Data" sheet:
#  - Row 2  : "Medstar POB North Tower" -> "Medstar POB South Tower" (Year Built 1967->1985, Gross Area 106517->76319)
#  - Row 3  : "1801 Pennsylvania Ave." -> "1801 Pennsylvania Avenue, LLC"
#  - Row 4  : Address "0300 E ST SW" -> "300 E ST SW", Owner "WELLS REIT/INDEPENDENCE SQUARE LLC" -> "TWO INDEPENDENCE HANA OW LLC", Gross Area 659773->627655
#  - Row 5  : Gross Area 65273->58717
#  - Row 6  : "President Madison Apartments" (a Multifamily property) replaced entirely with "Hampton House"
#             (Address, Postal Code, Year Built, Owner, Gross Area all updated)
#  - Row 7  : Postal Code 20007->20005, Gross Area 134036->145697
#  - Row 8  : Address "733 15TH ST NW" -> "1428 H ST NW"
#  - Row 9  : unchanged
#  - Row 10 : "DPW Vehicle Maintenance Facility 2" (a Service-Repair property) replaced entirely with
#             "School Without Walls @ Francis Stevens" (K-12 School)
#             (Property Type, Address, Postal Code, Year Built, Owner, Gross Area all updated)
#  - The "Year Built" column (I) loses its inherited date number-format (it now holds a plain year number).
#  - The "EUI Target Year" column (M) is removed entirely from the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("C2").Value = "Medstar POB South Tower"
$ws.Range("I2").Value = 1985
$ws.Range("L2").Value = 76319

# --- Row 3 ---
$ws.Range("C3").Value = "1801 Pennsylvania Avenue, LLC"

# --- Row 4 ---
$ws.Range("E4").Value = "300 E ST SW"
$ws.Range("J4").Value = "TWO INDEPENDENCE HANA OW LLC"
$ws.Range("L4").Value = 627655

# --- Row 5 ---
$ws.Range("L5").Value = 58717

# --- Row 6 ---
$ws.Range("C6").Value = "Hampton House"
$ws.Range("E6").Value = "2700 CONNECTICUT AVENUE NW"
$ws.Range("H6").Value = 20008
$ws.Range("I6").Value = 1921
$ws.Range("J6").Value = "2700 CONECTICUT AVENUE LLC"
$ws.Range("L6").Value = 83580

# --- Row 7 ---
$ws.Range("H7").Value = 20005
$ws.Range("L7").Value = 145697

# --- Row 8 ---
$ws.Range("E8").Value = "1428 H ST NW"

# --- Row 9 : unchanged ---

# --- Row 10 ---
$ws.Range("C10").Value = "School Without Walls @ Francis Stevens"
$ws.Range("D10").Value = "K-12 School"
$ws.Range("E10").Value = "2425 N STREET NW"
$ws.Range("H10").Value = 20037
$ws.Range("I10").Value = 1924
$ws.Range("J10").Value = "DISTRICT OF COLUMBIA"
$ws.Range("L10").Value = 127991

# The "Year Built" column no longer carries the inherited date format - reset to default/Normal style
$ws.Range("I2:I10").Style = "Normal"

# Remove the now-unused "EUI Target Year" column (M) entirely
[void]$ws.Range("M1:M1048576").Select()
$ws.Columns.Item(13).Delete()
